# Apply edit: replace shared-string placeholder ("-") cells with numeric 0,
# and flip sign on a handful of cells that were stored as negative values
# but should be positive (matching the target OOXML diff).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose shared-string text "-" placeholder becomes a literal numeric 0
$zeroCells = @(
    "G2","K2","L2","R2","V2","W2","X2","Y2","Z2","AC2","AD2","AE2","AF2","AV2","AW2","AX2","BD2",
    "BE2","BF2","BG2","BH2","BJ2","BK2","BP2","BQ2","BT2","CA2","CD2","CJ2","CO2","CP2","CR2","D3",
    "F3","G3","K3","L3","M3","O3","R3","V3","W3","X3","Y3","Z3","AC3","AD3","AE3","AF3","AV3","AW3",
    "AX3","BC3","BD3","BE3","BF3","BG3","BJ3","BN3","BO3","BP3","BQ3","BR3","BS3","BT3","BZ3","CA3",
    "CD3","CJ3","CO3","CP3","CR3","D4","F4","G4","K4","L4","M4","O4","R4","V4","W4","X4","Y4","Z4",
    "AC4","AD4","AE4","AF4","AV4","AW4","AX4","BC4","BD4","BE4","BF4","BG4","BJ4","BN4","BO4","BP4",
    "BQ4","BR4","BS4","BT4","BZ4","CA4","CD4","CJ4","CO4","CP4","CR4","F5","G5","K5","L5","M5","R5",
    "V5","W5","X5","Y5","Z5","AC5","AD5","AE5","AF5","AV5","AW5","AX5","BD5","BE5","BF5","BG5",
    "BJ5","BN5","BO5","BP5","BQ5","BR5","BS5","BT5","BX5","CD5","CJ5","CO5","CP5","CR5","D6","F6",
    "G6","K6","L6","M6","O6","R6","V6","W6","X6","Y6","Z6","AC6","AD6","AE6","AF6","AV6","AW6",
    "AX6","BC6","BD6","BE6","BF6","BG6","BJ6","BN6","BO6","BP6","BQ6","BR6","BS6","BT6","BZ6","CA6",
    "CD6","CJ6","CO6","CP6","CR6","S7","AI7","AJ7","AM7","AN7","AQ7","AR7","AS7","AT7","BI7","CH7",
    "CI7","CS7","S8","AI8","AJ8","AM8","AN8","AQ8","AR8","AS8","AT8","BI8","CH8","CI8","CS8","S9",
    "AI9","AJ9","AM9","AN9","AQ9","AR9","AS9","AT9","BI9","CH9","CI9","CS9","S10","AI10","AJ10",
    "AM10","AN10","AQ10","AR10","AS10","AT10","BI10","CH10","CI10","CS10","S11","AI11","AJ11",
    "AM11","AN11","AQ11","AR11","AS11","AT11","BI11","CH11","CI11","CS11","S12","AI12","AJ12",
    "AM12","AN12","AQ12","AR12","AS12","AT12","BI12","CH12","CI12","CS12","S13","AI13","AJ13",
    "AM13","AN13","AQ13","AR13","AS13","AT13","BI13","CH13","CI13","CS13","S14","AI14","AJ14",
    "AM14","AN14","AQ14","AR14","AS14","AT14","BI14","CH14","CI14","CS14","S15","AI15","AJ15",
    "AM15","AN15","AQ15","AR15","AS15","AT15","BI15","CH15","CI15","CS15","S16","AI16","AJ16",
    "AM16","AN16","AQ16","AR16","AS16","AT16","BI16","CH16","CI16","CS16"
)

foreach ($addr in $zeroCells) {
    $ws.Range($addr).Value = 0
}

# Cells that were negative and must become their positive (absolute) value
$positiveFixups = @{
    "CQ2" = 2192.21
    "CQ3" = 2621.72
    "CQ4" = 2452.5500000000002
    "CQ5" = 2122.62
    "CQ6" = 2205.2399999999998
    "CT6" = 812.84
    "CP7" = 195.48
    "CQ7" = 2793.51
    "CT7" = 566.51
    "CP8" = 631.62
    "CQ8" = 2161.89
    "CP9" = 1257.0899999999999
    "CQ9" = 904.79
    "CP10" = 1208.58
    "CP11" = 540.77
}

foreach ($addr in $positiveFixups.Keys) {
    $ws.Range($addr).Value = $positiveFixups[$addr]
}

# Restore the view state left by the author (selected cell W4)
$ws.Range("W4").Select()

